$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "270.43"
Set-TextValue $ws.Range("E2") "3.29%"
Set-TextValue $ws.Range("D3") "26.72"
Set-TextValue $ws.Range("E3") "-1.82%"
Set-TextValue $ws.Range("D4") "4.710"
Set-TextValue $ws.Range("E4") "0.15%"
Set-TextValue $ws.Range("D5") "0.06104"
Set-TextValue $ws.Range("E5") "-1.68%"
Set-TextValue $ws.Range("D6") "6.742"
Set-TextValue $ws.Range("E6") "0.33%"
Set-TextValue $ws.Range("D7") "0.8577"
Set-TextValue $ws.Range("E7") "0.84%"
Set-TextValue $ws.Range("D8") "0.8920"
Set-TextValue $ws.Range("D9") "0.1435"
Set-TextValue $ws.Range("E9") "1.68%"
Set-TextValue $ws.Range("D10") "0.04952"
Set-TextValue $ws.Range("E10") "6.90%"
Set-TextValue $ws.Range("D11") "0.07136"
Set-TextValue $ws.Range("E11") "0.70%"
Set-TextValue $ws.Range("D12") "0.03179"
Set-TextValue $ws.Range("E12") "-0.10%"
Set-TextValue $ws.Range("D13") "0.09035"
Set-TextValue $ws.Range("D14") "0.001529"
Set-TextValue $ws.Range("E14") "-0.89%"
Set-TextValue $ws.Range("D15") "0.0006084"
Set-TextValue $ws.Range("E15") "-1.16%"
Set-TextValue $ws.Range("D16") "0.005941"
Set-TextValue $ws.Range("E16") "-2.95%"
Set-TextValue $ws.Range("D17") "3.465"
Set-TextValue $ws.Range("E17") "-0.04%"
Set-TextValue $ws.Range("D18") "3.172"
Set-TextValue $ws.Range("E18") "0.16%"
Set-TextValue $ws.Range("E19") "2.93%"
Set-TextValue $ws.Range("D20") "0.3091"
Set-TextValue $ws.Range("E20") "-0.49%"
Set-TextValue $ws.Range("D21") "0.1309"
Set-TextValue $ws.Range("E21") "-0.12%"
Set-TextValue $ws.Range("D22") "3.842"
Set-TextValue $ws.Range("E22") "-5.86%"
Set-TextValue $ws.Range("D23") "0.04247"
Set-TextValue $ws.Range("E23") "0.09%"
Set-TextValue $ws.Range("D24") "0.001187"
Set-TextValue $ws.Range("E24") "-1.78%"
Set-TextValue $ws.Range("D25") "0.004151"
Set-TextValue $ws.Range("E25") "0.45%"
Set-TextValue $ws.Range("D26") "0.0001201"
Set-TextValue $ws.Range("E26") "0.02%"
Set-TextValue $ws.Range("E27") "5.02%"
Set-TextValue $ws.Range("D40") "0.03955"
Set-TextValue $ws.Range("E40") "1.19%"
Set-TextValue $ws.Range("D41") "0.1120"
Set-TextValue $ws.Range("E41") "0.64%"
Set-TextValue $ws.Range("D42") "0.004189"
Set-TextValue $ws.Range("E42") "1.43%"
Set-TextValue $ws.Range("D43") "0.002038"
Set-TextValue $ws.Range("E43") "-6.70%"
Set-TextValue $ws.Range("D44") "0.01277"
Set-TextValue $ws.Range("E44") "-8.27%"
Set-TextValue $ws.Range("D45") "0.00005122"
Set-TextValue $ws.Range("E45") "-0.97%"
Set-TextValue $ws.Range("D46") "0.00000000751"
Set-TextValue $ws.Range("E46") "0.01%"
Set-TextValue $ws.Range("D47") "0.9718"
Set-TextValue $ws.Range("E47") "479.87%"
Set-TextValue $ws.Range("D49") "0.00002102"
Set-TextValue $ws.Range("E49") "0.01%"
Set-TextValue $ws.Range("D50") "0.0002001"
Set-TextValue $ws.Range("E50") "0.01%"
